$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "dealerName"
$ws.Range("G2").Value = "dealerName1"
$ws.Range("G3").Value = "dealerName2"
$ws.Range("G4").Value = "dealerName3"
$ws.Range("G5").Value = "dealerName4"
$ws.Range("G6").Value = "dealerName5"
$ws.Range("G7").Value = "dealerName6"
$ws.Range("G8").Value = "dealerName7"
$ws.Range("G9").Value = "dealerName8"
$ws.Range("G10").Value = "dealerName9"
$ws.Range("G11").Value = "dealerName10"
$ws.Range("G12").Value = "dealerName11"
$ws.Range("G13").Value = "dealerName12"
$ws.Range("G14").Value = "dealerName13"

$ws.Range("G6").Select() | Out-Null
